$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.3106942176818848
$ws.Range("B1").Value = 0.3908534348011017
$ws.Range("C1").Value = 0.5684199929237366
$ws.Range("D1").Value = 2.31817626953125
$ws.Range("E1").Value = 5.575554847717285
